# Scenarij i tok događaja za rezervaciju sale.
#
# 1) The long description in B2 gets line breaks inserted (wrapped onto 3 lines).
# 2) A1:B10 (Naziv / Početna konfiguracija table) is turned into a real Excel
#    Table ("Table1", style TableStyleLight8, with header row + autofilter),
#    which is what drives the column-B width change, the row-2 height change
#    and the new A1:B12 selection left behind after the table was inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the wrapped description text in B2 -------------------------
$newText = "Na osnovu podataka koje unese administrator " + [char]10 + `
           "formira se baza podataka o prostorijama " + [char]10 + `
           "na fakultetu, te se kreiraju korisnički računi za osoblje"
$ws.Range("B2").Value = $newText
$ws.Range("B2").WrapText = $true

# Row grows to fit the now 3-line wrapped text.
$ws.Rows.Item(2).RowHeight = 42.75

# --- 2. Convert A1:B10 into an Excel Table ---------------------------------
$tableRange = $ws.Range("A1:B10")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight8"

# Column B narrows now that the long text wraps instead of running wide.
$ws.Columns.Item(2).ColumnWidth = 43.14

# --- 3. Leave the same selection behind as after inserting the table ------
$ws.Range("A1:B12").Select()
